$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.756.18'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.32%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.494.88'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.58%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.58'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.90%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.90'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.520.85'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.16%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0988'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.39%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.23%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.44'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.348'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.14%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.933.76'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.32%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.30'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.08%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.676.70'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.36%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000138'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.76%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.510.94'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.33%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.13'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.76%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.23'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.59%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '320.48'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.59%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.72'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.67'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.30%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.432'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -7.27%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.162'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.619.85'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.77%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.72%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.66'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.64'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.21%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0761'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.27%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.78'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.17%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.17'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -8.05%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.996'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.16%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.42'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.44%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.41'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.82%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.49'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.07%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.29'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.38%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.59'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.87%  '

$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.69'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.20%  '

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.54'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.34%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '296.19'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.14%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.61'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.91%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.804'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.81%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.995'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.26%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.600'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.76'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.57%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.46'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.23%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0923'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.92%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.50'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.16%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0226'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.62%  '
